$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 25 de Septiembre de 2020 a las 00:04"

# Row 4
$ws.Range("B4").Value = 7174692
$ws.Range("C4").Value = 35139
$ws.Range("D4").Value = 4419573
$ws.Range("E4").Value = 2547829
$ws.Range("G4").Value = 694
$ws.Range("H4").Value = 207290

# Row 6
$ws.Range("B6").Value = 4657702
$ws.Range("C6").Value = 29922
$ws.Range("D6").Value = 4023789
$ws.Range("E6").Value = 494105
$ws.Range("G6").Value = 743
$ws.Range("H6").Value = 139808

# Row 40
$ws.Range("B40").Value = 102513
$ws.Range("C40").Value = 138
$ws.Range("D40").Value = 92644
$ws.Range("E40").Value = 4034
$ws.Range("G40").Value = 13
$ws.Range("H40").Value = 5835

# Row 56
$ws.Range("B56").Value = 67701
$ws.Range("C56").Value = 687
$ws.Range("D56").Value = 60853
$ws.Range("E56").Value = 6617

# Row 61
$ws.Range("D61").Value = 42600
$ws.Range("E61").Value = 6831

# Row 84
$ws.Range("A84").Value = "Bulgaria"
$ws.Range("B84").Value = 19573
$ws.Range("C84").Value = 290
$ws.Range("D84").Value = 14013
$ws.Range("E84").Value = 4775
$ws.Range("G84").Value = 6
$ws.Range("H84").Value = 785

# Row 85
$ws.Range("A85").Value = "Costa de Marfil"
$ws.Range("B85").Value = 19501
$ws.Range("C85").Value = 71
$ws.Range("D85").Value = 19003
$ws.Range("E85").Value = 378
$ws.Range("H85").Value = 120

# Row 111
$ws.Range("B111").Value = 7433
$ws.Range("C111").Value = 8
$ws.Range("D111").Value = 7052
$ws.Range("E111").Value = 220

# Row 134
$ws.Range("B134").Value = 3966
$ws.Range("C134").Value = 42
$ws.Range("D134").Value = 1013
$ws.Range("E134").Value = 2770
$ws.Range("G134").Value = 2
$ws.Range("H134").Value = 183

# Row 145
$ws.Range("A145").Value = "Botsuana"
$ws.Range("B145").Value = 2921
$ws.Range("C145").Value = 354
$ws.Range("D145").Value = 701
$ws.Range("E145").Value = 2204
$ws.Range("G145").Value = 3
$ws.Range("H145").Value = 16

# Row 146
$ws.Range("A146").Value = "Malta"
$ws.Range("B146").Value = 2898
$ws.Range("C146").Value = 42
$ws.Range("D146").Value = 2191
$ws.Range("E146").Value = 680
$ws.Range("G146").Value = 2
$ws.Range("H146").Value = 27

# Row 147
$ws.Range("A147").Value = "Sudan del Sur"
$ws.Range("B147").Value = 2669
$ws.Range("C147").Value = 5
$ws.Range("D147").Value = 1290
$ws.Range("E147").Value = 1330
$ws.Range("H147").Value = 49

# Row 148
$ws.Range("B148").Value = 2579
$ws.Range("C148").Value = 44
$ws.Range("D148").Value = 1483
$ws.Range("E148").Value = 1025
$ws.Range("G148").Value = 2
$ws.Range("H148").Value = 71

# Row 166
$ws.Range("B166").Value = 1194
$ws.Range("C166").Value = 1
$ws.Range("E166").Value = 18

# Row 169
$ws.Range("B169").Value = 910
$ws.Range("C169").Value = 2
$ws.Range("D169").Value = 881
$ws.Range("E169").Value = 14

# Row 215
$ws.Range("A215").Value = "Montserrat"
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1

# Row 216
$ws.Range("A216").Value = "Islas Malvinas"
$ws.Range("D216").Value = 13
$ws.Range("H216").Value = 0
